# Convert the ingredient weight column (N) from "jin"/"顆" units (e.g. "12斤",
# "8顆") to plain decimal-kg text values (e.g. "7.199999999999999").
#
# A direct `$ws.Range("N2").Value = "15.0"` would make Excel re-interpret the
# numeric-looking string as a *number*, losing the original text cell type.
# Instead we build each value as formula-derived TEXT in a scratch cell
# (`="15.0"` -> the string "15.0"), copy it, and paste *values only* into the
# destination. That preserves the destination's text type/format exactly as
# it already was, with no NumberFormat/style changes needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$xlPasteValues = -4163

function Set-TextValue($cellAddress, $text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellAddress).PasteSpecial($xlPasteValues)
}

Set-TextValue "N2"  "15.0"
Set-TextValue "N3"  "7.199999999999999"
Set-TextValue "N4"  "1.7999999999999998"
Set-TextValue "N5"  "1.7999999999999998"
Set-TextValue "N6"  "4.8"
Set-TextValue "N7"  "2.4"
Set-TextValue "N8"  "0.6"
Set-TextValue "N9"  "0.6"
Set-TextValue "N10" "0.6"

# Remove the scratch cell so it doesn't linger in the saved workbook.
$scratch.ClearContents()
